$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "57.007.17"
$ws.Range("E2").Value = "  +11.21%  "
$ws.Range("D3").Value = "3.269.71"
$ws.Range("E3").Value = "  +6.79%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.66%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0973"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.83%  "
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "3.779.66"
$ws.Range("E13").Value = "  +6.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").Value = "3.263.13"
$ws.Range("E16").Value = "  +6.47%  "
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "56.871.69"
$ws.Range("E19").Value = "  +10.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.30%  "
$ws.Range("E21").Value = "  +11.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "306.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.66%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +5.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0482"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +24.22%  "
$ws.Range("E39").Value = "  +7.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "134.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.36%  "
$ws.Range("E42").Value = "  +5.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("E44").Value = "  +4.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "2.158.18"
$ws.Range("E48").Value = "  +4.56%  "
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("E50").Value = "  +43.27%  "
$ws.Range("E51").Value = "  -4.49%  "
